$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "No. CHN-" marking number
$ws.Range("C2").Value = "250215W99-BA"

# "Marking : " value
$ws.Range("C7").Value = "279-W99-BA-SEA"

# First cargo line (row 11): Koli qty, item name, weight qty, weight (M3)
$ws.Range("D11").Value = 4
$ws.Range("G11").Value = "S0083449  Mainan"
$ws.Range("M11").Value = 37
$ws.Range("N11").Value = 0.5055

# Second cargo line (row 12): Koli qty, unit label, item name, weight qty, weight (M3)
$ws.Range("D12").Value = 3
$ws.Range("E12").Value = "Koli"
$ws.Range("G12").Value = "S0083517  Mainan"
$ws.Range("M12").Value = 20
$ws.Range("N12").Value = 0.2061

$wb.Save()
